# Adds new Array problems to the Notes workbook:
#  - Adds a missing hyperlink to the existing "Product of Array Except Self" row (D11)
#  - Adds a new row for "Maximum Subarray"
#  - Adds a new row for "Find Minimum in Rotated Sorted Array"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: new "Maximum Subarray" row ---
# Copy the formatting used by the other "Good"-styled name cells (e.g. A2) so the
# new cell picks up the exact same cell style instead of a brand new one.
$ws.Range("A2").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = "Maximum Subarray"
$ws.Range("B12").Value = "Return largest sum of the sub-array"
$ws.Range("C12").Value = "Take maximum of current sum + current num and current num to keep current sum positive. Keep updating maximum sum using current sum"
$ws.Hyperlinks.Add($ws.Range("D12"), "https://leetcode.com/problems/product-of-array-except-self/")
$ws.Range("D12").Style = "Hyperlink"

# --- Row 11: existing "Product of Array Except Self" row was missing its Link cell ---
$ws.Hyperlinks.Add($ws.Range("D11"), "https://leetcode.com/problems/product-of-array-except-self/")
$ws.Range("D11").Style = "Hyperlink"

# --- Row 13: new "Find Minimum in Rotated Sorted Array" row ---
# Copy the formatting used by the other "Neutral"-styled name cells (e.g. A10).
$ws.Range("A10").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "Find Minimum in Rotated Sorted Array"
$ws.Range("B13").Value = "Return min in O(log n) time"
$ws.Range("C13").Value = "Use modified binary search. Compare mid element with first and last element. If n[mid] > n[right], element is in right. If n[mid] < n[left], element is in left. Else return n[left]"
$ws.Hyperlinks.Add($ws.Range("D13"), "https://leetcode.com/problems/find-minimum-in-rotated-sorted-array/")
$ws.Range("D13").Style = "Hyperlink"

# --- Misc view state tweaks captured in the diff ---
[void]$ws.Range("C4").Select()

Write-Host "Added new Array problems"
